$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.25
$ws.Range("C2").Value = 0.4375
$ws.Range("J2").Value = 0.03125
$ws.Range("P2").Value = 0.25
$ws.Range("S2").Value = 0.03125
$ws.Range("C3").Value = 0.06666666666666667
$ws.Range("J3").Value = 0.06666666666666667
$ws.Range("P3").Value = 0.5333333333333333
$ws.Range("S3").Value = 0.3333333333333333
$ws.Range("P4").Value = 0.75
$ws.Range("S4").Value = 0.25
$ws.Range("J6").Value = 0.3076923076923077
$ws.Range("Q6").Value = 0.3076923076923077
$ws.Range("R6").Value = 0.1538461538461539
$ws.Range("S6").Value = 0.2307692307692308
$ws.Range("B7").Value = 0.1111111111111111
$ws.Range("J7").Value = 0.3333333333333333
$ws.Range("Q7").Value = 0.2222222222222222
$ws.Range("R7").Value = 0.1111111111111111
$ws.Range("S7").Value = 0.2222222222222222
$ws.Range("B8").Value = 0.07692307692307693
$ws.Range("F8").Value = 0.03846153846153846
$ws.Range("J8").Value = 0.1538461538461539
$ws.Range("Q8").Value = 0.07692307692307693
$ws.Range("R8").Value = 0.1153846153846154
$ws.Range("S8").Value = 0.5384615384615384
$ws.Range("B9").Value = 0.1
$ws.Range("F9").Value = 0.3
$ws.Range("J9").Value = 0.2
$ws.Range("Q9").Value = 0.2
$ws.Range("R9").Value = 0.1
$ws.Range("S9").Value = 0.1
$ws.Range("B10").Value = 0.1739130434782609
$ws.Range("D10").Value = 0.03478260869565217
$ws.Range("F10").Value = 0.04347826086956522
$ws.Range("J10").Value = 0.1043478260869565
$ws.Range("O10").Value = 0.008695652173913044
$ws.Range("Q10").Value = 0.2956521739130435
$ws.Range("R10").Value = 0.08695652173913043
$ws.Range("S10").Value = 0.2521739130434782
$ws.Range("G11").Value = 0.1818181818181818
$ws.Range("K11").Value = 0.09090909090909091
$ws.Range("L11").Value = 0.7272727272727273
$ws.Range("G12").Value = 0.875
$ws.Range("J12").Value = 0.125
$ws.Range("J13").Value = 1
$ws.Range("H15").Value = 0.06666666666666667
$ws.Range("J15").Value = 0.5333333333333333
$ws.Range("K15").Value = 0.06666666666666667
$ws.Range("M15").Value = 0.06666666666666667
$ws.Range("O15").Value = 0.1333333333333333
$ws.Range("S15").Value = 0.1333333333333333
$ws.Range("H16").Value = 0.1578947368421053
$ws.Range("J16").Value = 0.631578947368421
$ws.Range("K16").Value = 0.05263157894736842
$ws.Range("O16").Value = 0.05263157894736842
$ws.Range("S16").Value = 0.1052631578947368
$ws.Range("H17").Value = 0.09090909090909091
$ws.Range("I17").Value = 0.09090909090909091
$ws.Range("J17").Value = 0.5681818181818182
$ws.Range("K17").Value = 0.06818181818181818
$ws.Range("O17").Value = 0.04545454545454546
$ws.Range("S17").Value = 0.1363636363636364
$ws.Range("H18").Value = 0.25
$ws.Range("J18").Value = 0.5
$ws.Range("K18").Value = 0.0625
$ws.Range("S18").Value = 0.1875
$ws.Range("F19").Value = 0.01388888888888889
$ws.Range("H19").Value = 0.2083333333333333
$ws.Range("I19").Value = 0.09722222222222222
$ws.Range("J19").Value = 0.4722222222222222
$ws.Range("K19").Value = 0.04166666666666666
$ws.Range("O19").Value = 0.1111111111111111
$ws.Range("S19").Value = 0.05555555555555555
